{"js": "// 1. \"Yes\" -> \"OK\" (the bold confirmation answer to the \"program the\n//    configuration memory now?\" prompt).\nconst yesResults = context.document.body.search(\"Yes\", { matchCase: true, matchWholeWord: true });\nyesResults.load(\"text\");\nawait context.sync();\nif (yesResults.items.length > 0) {\n  yesResults.items[0].insertText(\"OK\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. \"saturnprom\" -> \"saturn\" + \"fallback\" (now reads \"saturnfallback\").\nconst saturnResults = context.document.body.search(\"saturnprom\", { matchCase: true });\nsaturnResults.load(\"text\");\nawait context.sync();\nif (saturnResults.items.length > 0) {\n  const saturnRange = saturnResults.items[0];\n  saturnRange.insertText(\"saturn\", Word.InsertLocation.replace);\n  await context.sync();\n\n  const fallbackAnchor = context.document.body.search(\"saturn\", { matchCase: true, matchWholeWord: true });\n  fallbackAnchor.load(\"text\");\n  await context.sync();\n  fallbackAnchor.items[0].insertText(\"fallback\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 3. Insert a new list paragraph \"Click OK to begin\" right before the\n//    \"The config prom is programmed...\" list paragraph.\nconst targetResults = context.document.body.search(\"The config prom is programmed\", { matchCase: true });\ntargetResults.load(\"text\");\nawait context.sync();\nif (targetResults.items.length > 0) {\n  const targetParagraph = targetResults.items[0].paragraphs.getFirst();\n  const newParagraph = targetParagraph.insertParagraph(\"Click OK to begin\", Word.InsertLocation.before);\n  newParagraph.styleBuiltIn = Word.Style.listParagraph;\n  newParagraph.attachToList(2, 0);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"Yes\" -> \"OK\" (the bold confirmation answer to the \"program the\n#    configuration memory now?\" prompt).\n$find = $d.Content.Find\n$find.ClearFormatting()\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n$find.Execute(\"Yes\", $true, $true, $false, $false, $false, $true, 1, $false, \"OK\", 2) | Out-Null\n\n# 2. \"saturnprom\" -> \"saturn\" + \"fallback\" (now reads \"saturnfallback\").\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"saturnprom\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $saturnRange = $find2.Duplicate\n    $saturnRange.Text = \"saturn\"\n    $saturnRange.Collapse(0)\n    $saturnRange.InsertAfter(\"fallback\")\n}\n\n# 3. Insert a new list paragraph \"Click OK to begin\" right before the\n#    \"The config prom is programmed...\" list paragraph.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*The config prom is programmed*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $insertPoint = $target.Range.Duplicate\n    $insertPoint.Collapse(1)\n    $insertPoint.InsertParagraphBefore()\n\n    $newPara = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*The config prom is programmed*\") {\n            break\n        }\n        $newPara = $p\n    }\n    $newPara.Range.Text = \"Click OK to begin\"\n}\n\nWrite-Output \"done\"\n"}
